# TC21: add a "StatQuery" column between the existing "query" and "dbExcel"
# columns, holding the Neo4j stats query used to produce file/case/trial counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; existing B ("dbExcel") and C ("WebExcel") shift
# right to become C and D. Excel's default insert behavior formats the new
# column like the one to its left (column A), which carries row 2's
# wrap-text style - matching the target's A2/B2 "s=1" formatting.
$ws.Columns("B").Insert()

# Header + query text for the new "StatQuery" column.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Prostate cancer, NOS']  OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Column B should be as wide as column A, and row 2 should wrap its text
# (same look as the original query column).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Range("B2").WrapText = $true

# New active selection, matching the workbook's saved cursor position.
$ws.Range("A2").Select()
